$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("E2").Value = 11.6783744
$ws.Range("E3").Value = 2.1624928
$ws.Range("E4").Value = 0.9023760000000001
$ws.Range("E5").Value = 1.430928
$ws.Range("E6").Value = 0.9591919999999999
$ws.Range("E7").Value = 0.2855808
$ws.Range("E8").Value = 0.0899712
$ws.Range("C9").Value = 837
$ws.Range("E9").Value = 0.8639312
$ws.Range("C10").Value = 601
$ws.Range("E10").Value = 0.5197824000000001
$ws.Range("C11").Value = 1795
$ws.Range("E11").Value = 0.753816
$ws.Range("C12").Value = 412
$ws.Range("E12").Value = 1.2334432
$ws.Range("C13").Value = 3247
$ws.Range("E13").Value = 12.55878399999999
$ws.Range("C14").Value = 1192
$ws.Range("E14").Value = 4.825121600000001
$ws.Range("C15").Value = 302
$ws.Range("E15").Value = 0.8808624
$ws.Range("C16").Value = 419
$ws.Range("E16").Value = 0.8803135999999999
$ws.Range("C17").Value = 596
$ws.Range("E17").Value = 1.39932
$ws.Range("C18").Value = 173
$ws.Range("E18").Value = 0.447552
$ws.Range("C19").Value = 61
$ws.Range("E19").Value = 0.1245056
$ws.Range("E20").Value = 62.88014400000002
$ws.Range("E21").Value = 66.75336960000001
$ws.Range("E22").Value = 82.62881280000002
$ws.Range("E23").Value = 254.2461792000001

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("C2").Value = 595
$ws.Range("D2").Value = 145
$ws.Range("F2").Value = 306
$ws.Range("G2").Value = 97
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 419
$ws.Range("J2").Value = 299
$ws.Range("K2").Value = 896
$ws.Range("L2").Value = 202
$ws.Range("C3").Value = 600
$ws.Range("D3").Value = 150
$ws.Range("F3").Value = 293
$ws.Range("G3").Value = 85
$ws.Range("H3").Value = 27
$ws.Range("I3").Value = 418
$ws.Range("K3").Value = 899
$ws.Range("L3").Value = 210
$ws.Range("C4").Value = 595
$ws.Range("D4").Value = 153
$ws.Range("F4").Value = 294
$ws.Range("G4").Value = 80
$ws.Range("H4").Value = 26
$ws.Range("I4").Value = 420
$ws.Range("K4").Value = 905
$ws.Range("L4").Value = 217
$ws.Range("C5").Value = 597
$ws.Range("D5").Value = 149
$ws.Range("F5").Value = 302
$ws.Range("G5").Value = 93
$ws.Range("H5").Value = 35
$ws.Range("J5").Value = 303
$ws.Range("K5").Value = 897
$ws.Range("L5").Value = 202
$ws.Range("C6").Value = 590
$ws.Range("D6").Value = 150
$ws.Range("F6").Value = 300
$ws.Range("G6").Value = 89
$ws.Range("H6").Value = 23
$ws.Range("I6").Value = 423
$ws.Range("J6").Value = 299
$ws.Range("K6").Value = 890
$ws.Range("L6").Value = 215

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("C2").Value = 74.375
$ws.Range("D2").Value = 18.125
$ws.Range("F2").Value = 38.25
$ws.Range("G2").Value = 12.125
$ws.Range("H2").Value = 3.875
$ws.Range("I2").Value = 52.375
$ws.Range("J2").Value = 37.375
$ws.Range("K2").Value = 112
$ws.Range("L2").Value = 25.25
$ws.Range("C3").Value = 112.5
$ws.Range("D3").Value = 28.125
$ws.Range("F3").Value = 54.9375
$ws.Range("G3").Value = 15.9375
$ws.Range("H3").Value = 5.0625
$ws.Range("I3").Value = 78.375
$ws.Range("K3").Value = 168.5625
$ws.Range("L3").Value = 39.375
$ws.Range("C4").Value = 130.15625
$ws.Range("D4").Value = 33.46875
$ws.Range("F4").Value = 64.3125
$ws.Range("G4").Value = 17.5
$ws.Range("H4").Value = 5.6875
$ws.Range("I4").Value = 91.875
$ws.Range("K4").Value = 197.96875
$ws.Range("L4").Value = 47.46875
$ws.Range("C5").Value = 139.921875
$ws.Range("D5").Value = 34.921875
$ws.Range("F5").Value = 70.78125
$ws.Range("G5").Value = 21.796875
$ws.Range("H5").Value = 8.203125
$ws.Range("J5").Value = 71.015625
$ws.Range("K5").Value = 210.234375
$ws.Range("L5").Value = 47.34375
$ws.Range("C6").Value = 142.890625
$ws.Range("D6").Value = 36.328125
$ws.Range("F6").Value = 72.65625
$ws.Range("G6").Value = 21.5546875
$ws.Range("H6").Value = 5.5703125
$ws.Range("I6").Value = 102.4453125
$ws.Range("J6").Value = 72.4140625
$ws.Range("K6").Value = 215.546875
$ws.Range("L6").Value = 52.0703125

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 50470
$ws.Range("B3").Value = 523005
$ws.Range("B4").Value = 523005
$ws.Range("B5").Value = 4931190.000000001
